# Language.xlsx edit: change the "reset" button instruction/confirmation flow.
# - Replace the old "reset_instruction" row with a new "clear_objects_confirm"
#   confirmation prompt.
# - Add new "yes"/"no" localization keys (with YES/NO display values) used by
#   the new confirmation dialog buttons.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the (soon to be repurposed) row 50 so the
# new yes/no keys land at rows 51-52, pushing everything below down by two.
$ws.Rows("51:52").Insert()

# Row 50: repurpose reset_instruction -> clear_objects_confirm
$ws.Range("A50").Value = "clear_objects_confirm"

# Row 51: new "yes" key
$ws.Range("A51").Value = "yes"

# Row 52: new "no" key
$ws.Range("A52").Value = "no"

# Fill in the value column after the keys, matching the order new strings
# were introduced in the authored workbook.
$ws.Range("B51").Value = "YES"
$ws.Range("B52").Value = "NO"
$ws.Range("B50").Value = "Do you want to remove all the objects placed in the world?"

# Update the view so the active cell matches where the edit took place.
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A51").Select()
